$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 101
$ws.Range("H101").Value = 1646.25
$ws.Range("I101").Value = 1712.1428
$ws.Range("K101").Value = 5136.428400000001
$ws.Range("M101").Value = -3514.428400000001
# row 138
$ws.Range("H138").Value = 1810.8298
$ws.Range("I138").Value = 1159.7428
$ws.Range("J138").Value = 3709.8333
$ws.Range("K138").Value = 3479.2284
$ws.Range("L138").Value = 11129.4999
$ws.Range("M138").Value = 1660.7716
$ws.Range("N138").Value = -21409.4999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 17216.988
$ws.Range("I32").Value = 3651.378
$ws.Range("K32").Value = 3651.378
$ws.Range("M32").Value = -3364.378
# row 45
$ws.Range("H45").Value = 3591.75
$ws.Range("I45").Value = 2330.5625
$ws.Range("J45").Value = 6114.125
$ws.Range("K45").Value = 2330.5625
$ws.Range("L45").Value = 6114.125
$ws.Range("M45").Value = -1953.5625
$ws.Range("N45").Value = -6868.125
# row 61
$ws.Range("H61").Value = 2052.5625
$ws.Range("I61").Value = 2101.7856
$ws.Range("J61").Value = 1708
$ws.Range("K61").Value = 2101.7856
$ws.Range("L61").Value = 1708
$ws.Range("M61").Value = -1889.7856
$ws.Range("N61").Value = -2132
# row 122
$ws.Range("H122").Value = 2168.1333
$ws.Range("I122").Value = 2178.6155
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 6535.8465
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -4085.8465
$ws.Range("N122").Value = -11200
# row 136
$ws.Range("H136").Value = 2052.5625
$ws.Range("I136").Value = 2101.7856
$ws.Range("J136").Value = 1708
$ws.Range("K136").Value = 6305.3568
$ws.Range("L136").Value = 5124
$ws.Range("M136").Value = -3755.3568
$ws.Range("N136").Value = -10224

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 99
$ws.Range("H99").Value = 1971.2307
$ws.Range("I99").Value = 1826.5
$ws.Range("J99").Value = 2202.8
$ws.Range("K99").Value = 1826.5
$ws.Range("L99").Value = 2202.8
$ws.Range("M99").Value = -328.5
$ws.Range("N99").Value = -5198.8
# row 122
$ws.Range("H122").Value = 1465.7826
$ws.Range("I122").Value = 965.4706
$ws.Range("J122").Value = 2883.3333
$ws.Range("K122").Value = 2896.4118
$ws.Range("L122").Value = 8649.999899999999
$ws.Range("M122").Value = -446.4117999999999
$ws.Range("N122").Value = -13549.9999
# row 126
$ws.Range("H126").Value = 1971.2307
$ws.Range("I126").Value = 1826.5
$ws.Range("J126").Value = 2202.8
$ws.Range("K126").Value = 5479.5
$ws.Range("L126").Value = 6608.400000000001
$ws.Range("M126").Value = -3009.5
$ws.Range("N126").Value = -11548.4
# row 132
$ws.Range("H132").Value = 3272.4849
$ws.Range("I132").Value = 2692.0386
$ws.Range("J132").Value = 5428.4287
$ws.Range("K132").Value = 8076.1158
$ws.Range("L132").Value = 16285.2861
$ws.Range("M132").Value = -5546.1158
$ws.Range("N132").Value = -21345.2861
# row 134
$ws.Range("H134").Value = 5268942
$ws.Range("I134").Value = 7286.857
$ws.Range("J134").Value = 20001576
$ws.Range("K134").Value = 21860.571
$ws.Range("L134").Value = 60004728
$ws.Range("M134").Value = -19325.571
$ws.Range("N134").Value = -60009798

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 785063.0600000001
$ws.Range("I5").Value = 719.5
$ws.Range("J5").Value = 2667487.5
$ws.Range("K5").Value = 2158.5
$ws.Range("L5").Value = 8002462.5
$ws.Range("M5").Value = -2046.5
$ws.Range("N5").Value = -8002686.5
# row 80
$ws.Range("H80").Value = 6373.3335
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 6542.857
$ws.Range("K80").Value = 12000
$ws.Range("L80").Value = 19628.571
$ws.Range("M80").Value = -11064
$ws.Range("N80").Value = -21500.571
# row 83
$ws.Range("H83").Value = 6373.3335
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 6542.857
$ws.Range("K83").Value = 36000
$ws.Range("L83").Value = 58885.713
$ws.Range("M83").Value = -31320
$ws.Range("N83").Value = -68245.713
# row 112
$ws.Range("H112").Value = 48654736
$ws.Range("I112").Value = 1142.3334
$ws.Range("J112").Value = 55605250
$ws.Range("K112").Value = 3427.0002
$ws.Range("L112").Value = 166815750
$ws.Range("M112").Value = -2319.0002
$ws.Range("N112").Value = -166817966
# row 122
$ws.Range("H122").Value = 46528.773
$ws.Range("J122").Value = 52164.145
$ws.Range("L122").Value = 469477.305
$ws.Range("N122").Value = -474377.305
# row 135
$ws.Range("H135").Value = 785063.0600000001
$ws.Range("I135").Value = 719.5
$ws.Range("J135").Value = 2667487.5
$ws.Range("K135").Value = 6475.5
$ws.Range("L135").Value = 24007387.5
$ws.Range("M135").Value = -3940.5
$ws.Range("N135").Value = -24012457.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 126
$ws.Range("H126").Value = 10935.071
$ws.Range("I126").Value = 2120.3635
$ws.Range("J126").Value = 16638.705
$ws.Range("K126").Value = 6361.0905
$ws.Range("L126").Value = 49916.11500000001
$ws.Range("M126").Value = -3891.0905
$ws.Range("N126").Value = -54856.11500000001
# row 132
$ws.Range("H132").Value = 3679.7646
$ws.Range("I132").Value = 3638.5454
$ws.Range("J132").Value = 3755.3333
$ws.Range("K132").Value = 10915.6362
$ws.Range("L132").Value = 11265.9999
$ws.Range("M132").Value = -8385.636200000001
$ws.Range("N132").Value = -16325.9999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 9093464
$ws.Range("I7").Value = 14287671
$ws.Range("J7").Value = 3599.75
$ws.Range("K7").Value = 14287671
$ws.Range("L7").Value = 3599.75
$ws.Range("M7").Value = -14287559
$ws.Range("N7").Value = -3823.75
# row 40
$ws.Range("H40").Value = 2068.6924
$ws.Range("I40").Value = 1799.3
$ws.Range("J40").Value = 2966.6667
$ws.Range("K40").Value = 1799.3
$ws.Range("L40").Value = 2966.6667
$ws.Range("M40").Value = -1663.3
$ws.Range("N40").Value = -3238.6667
# row 122
$ws.Range("H122").Value = 3208.2354
$ws.Range("I122").Value = 2375
$ws.Range("J122").Value = 3948.889
$ws.Range("K122").Value = 7125
$ws.Range("L122").Value = 11846.667
$ws.Range("M122").Value = -4675
$ws.Range("N122").Value = -16746.667
# row 126
$ws.Range("H126").Value = 9093464
$ws.Range("I126").Value = 14287671
$ws.Range("J126").Value = 3599.75
$ws.Range("K126").Value = 42863013
$ws.Range("L126").Value = 10799.25
$ws.Range("M126").Value = -42860543
$ws.Range("N126").Value = -15739.25
# row 132
$ws.Range("H132").Value = 2058.535
$ws.Range("I132").Value = 1760.5428
$ws.Range("J132").Value = 3362.25
$ws.Range("K132").Value = 5281.6284
$ws.Range("L132").Value = 10086.75
$ws.Range("M132").Value = -2751.6284
$ws.Range("N132").Value = -15146.75

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 2288.8445
$ws.Range("I132").Value = 2287.3076
$ws.Range("J132").Value = 2298.8333
$ws.Range("K132").Value = 6861.9228
$ws.Range("L132").Value = 6896.499899999999
$ws.Range("M132").Value = -4331.9228
$ws.Range("N132").Value = -11956.4999
